# "feat: added weekend only feature"
#
# Each of the 6 fare sheets (K2, H1, H2, P, L, K1) gets two new trailing
# columns F ("Unnamed: 6") and G ("Unnamed: 7") added by the re-export
# (pandas df.to_excel of a frame that picked up two extra blank columns).
# On the "direct" sheets (H1, P, L, K1) the booking-class / fare data also
# changed, and the P / L sheets lost their extra SGN/BKK rows (only the
# TPE rows remain). K2 and H2 are untouched other than the already-empty
# "direct" column cells disappearing (Excel drops wholly-empty cells with
# no sibling data on save).

$wb = $excel.ActiveWorkbook

function Add-UnnamedHeaders($ws) {
    # Copy the format of the existing "direct" header (column E) onto the
    # two new trailing header cells, then set their text.
    $ws.Cells.Item(1, 5).Copy()
    $ws.Range("F1:G1").PasteSpecial(-4122)
    $ws.Cells.Item(1, 6).Value = "Unnamed: 6"
    $ws.Cells.Item(1, 7).Value = "Unnamed: 7"
}

function Set-BlankCell($ws, $row, $col) {
    # Force a genuinely empty *text* cell to materialise (mirrors the
    # empty inlineStr cells pandas/openpyxl emit for NaN string columns).
    $ws.Cells.Item($row, $col).Value = "'"
}

# ---------------------------------------------------------------------
# K2 (sheet1) — no data/header changes; the empty "direct" cells (E2:E4)
# just disappear on save since they carry no content.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("K2")
$ws.Cells.Item(2, 5).ClearContents()
$ws.Cells.Item(3, 5).ClearContents()
$ws.Cells.Item(4, 5).ClearContents()

# ---------------------------------------------------------------------
# H1 (sheet2)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("H1")
Add-UnnamedHeaders $ws

$ws.Cells.Item(2, 2).Value = "J"
$ws.Cells.Item(2, 4).Value = 9765
Set-BlankCell $ws 2 6
Set-BlankCell $ws 2 7

$ws.Cells.Item(3, 2).Value = "C"
$ws.Cells.Item(3, 4).Value = 7800
Set-BlankCell $ws 3 6
Set-BlankCell $ws 3 7

$ws.Cells.Item(4, 1).Value = "TPE"
$ws.Cells.Item(4, 2).Value = "D"
$ws.Cells.Item(4, 4).Value = 5700
Set-BlankCell $ws 4 6
Set-BlankCell $ws 4 7

# ---------------------------------------------------------------------
# H2 (sheet3) — no data/header changes, same as K2.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("H2")
$ws.Cells.Item(2, 5).ClearContents()
$ws.Cells.Item(3, 5).ClearContents()
$ws.Cells.Item(4, 5).ClearContents()

# ---------------------------------------------------------------------
# P (sheet4) — drop rows 5 & 6 (old SGN/BKK rows), keep only the 3 TPE
# rows, add the two trailing columns.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("P")
$ws.Rows.Item(5).Delete()
$ws.Rows.Item(5).Delete()

Add-UnnamedHeaders $ws

$ws.Cells.Item(2, 2).Value = "J"
$ws.Cells.Item(2, 4).Value = 10120
Set-BlankCell $ws 2 6
Set-BlankCell $ws 2 7

$ws.Cells.Item(3, 2).Value = "C"
$ws.Cells.Item(3, 4).Value = 8030
Set-BlankCell $ws 3 6
Set-BlankCell $ws 3 7

$ws.Cells.Item(4, 1).Value = "TPE"
$ws.Cells.Item(4, 2).Value = "D"
$ws.Cells.Item(4, 4).Value = 5900
Set-BlankCell $ws 4 6
Set-BlankCell $ws 4 7

# ---------------------------------------------------------------------
# L (sheet5) — drop rows 5, 6 & 7, keep only the 3 TPE rows, add the two
# trailing columns; G2 gets a real remarks note instead of staying blank.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("L")
$ws.Rows.Item(5).Delete()
$ws.Rows.Item(5).Delete()
$ws.Rows.Item(5).Delete()

Add-UnnamedHeaders $ws

$ws.Cells.Item(2, 2).Value = "J"
$ws.Cells.Item(2, 4).Value = 9275
Set-BlankCell $ws 2 6
$ws.Cells.Item(2, 7).Value = "Remarks:`n'O' for One Season"

$ws.Cells.Item(3, 2).Value = "C"
$ws.Cells.Item(3, 4).Value = 7430
Set-BlankCell $ws 3 6
Set-BlankCell $ws 3 7

$ws.Cells.Item(4, 2).Value = "D"
$ws.Cells.Item(4, 4).Value = 5400
Set-BlankCell $ws 4 6
Set-BlankCell $ws 4 7

# ---------------------------------------------------------------------
# K1 (sheet6)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("K1")
Add-UnnamedHeaders $ws

$ws.Cells.Item(2, 2).Value = "J"
$ws.Cells.Item(2, 4).Value = 9655
Set-BlankCell $ws 2 6
Set-BlankCell $ws 2 7

$ws.Cells.Item(3, 2).Value = "C"
$ws.Cells.Item(3, 4).Value = 7720
Set-BlankCell $ws 3 6
Set-BlankCell $ws 3 7

$ws.Cells.Item(4, 1).Value = "TPE"
$ws.Cells.Item(4, 2).Value = "D"
$ws.Cells.Item(4, 4).Value = 5630
Set-BlankCell $ws 4 6
Set-BlankCell $ws 4 7
